$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.402.25'
$ws.Range('E2').Value = '  -1.30%  '
$ws.Range('D3').Value = '3.326.69'
$ws.Range('E3').Value = '  -3.22%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''247.91'
$ws.Range('E5').Value = '  -4.48%  '
$ws.Range('D6').Value = '''652.75'
$ws.Range('E6').Value = '  -2.75%  '
$ws.Range('D7').Value = '''1.38'
$ws.Range('E7').Value = '  -11.10%  '
$ws.Range('D8').Value = '''0.419'
$ws.Range('E8').Value = '  -10.53%  '
$ws.Range('D10').Value = '''1.00'
$ws.Range('E10').Value = '  -8.47%  '
$ws.Range('D11').Value = '3.322.90'
$ws.Range('E11').Value = '  -3.24%  '
$ws.Range('E12').Value = '  -6.93%  '
$ws.Range('D13').Value = '''40.22'
$ws.Range('E13').Value = '  -6.83%  '
$ws.Range('D14').Value = '97.285.15'
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('D15').Value = '''6.03'
$ws.Range('E15').Value = '  -1.84%  '
$ws.Range('D16').Value = '''0.0000252'
$ws.Range('E16').Value = '  -8.57%  '
$ws.Range('D17').Value = '3.948.23'
$ws.Range('E17').Value = '  -2.98%  '
$ws.Range('D18').Value = '''8.65'
$ws.Range('E18').Value = '  +6.39%  '
$ws.Range('D19').Value = '3.326.86'
$ws.Range('E19').Value = '  -3.12%  '
$ws.Range('D20').Value = '''0.523'
$ws.Range('E20').Value = '  +19.75%  '
$ws.Range('D21').Value = '''16.73'
$ws.Range('E21').Value = '  -3.87%  '
$ws.Range('D22').Value = '''10.54'
$ws.Range('E22').Value = '  -1.69%  '
$ws.Range('D23').Value = '''495.40'
$ws.Range('E23').Value = '  -7.89%  '
$ws.Range('D24').Value = '''3.29'
$ws.Range('E24').Value = '  -8.74%  '
$ws.Range('E25').Value = '  -9.55%  '
$ws.Range('D26').Value = '''6.47'
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('D27').Value = '''95.58'
$ws.Range('E27').Value = '  -6.90%  '
$ws.Range('D28').Value = '''12.01'
$ws.Range('E28').Value = '  -7.15%  '
$ws.Range('D29').Value = '3.503.41'
$ws.Range('E29').Value = '  -2.85%  '
$ws.Range('E30').Value = '  -4.10%  '
$ws.Range('D31').Value = '''0.993'
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('D32').Value = '''10.84'
$ws.Range('E32').Value = '  -7.10%  '
$ws.Range('D33').Value = '''0.187'
$ws.Range('E33').Value = '  -6.63%  '
$ws.Range('E34').Value = '  +11.04%  '
$ws.Range('D35').Value = '''0.997'
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('D36').Value = '''0.548'
$ws.Range('E36').Value = '  -2.92%  '
$ws.Range('D37').Value = '''28.21'
$ws.Range('E37').Value = '  -7.44%  '
$ws.Range('D38').Value = '''1.46'
$ws.Range('E38').Value = '  +3.21%  '
$ws.Range('D39').Value = '''7.54'
$ws.Range('E39').Value = '  -5.16%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('E41').Value = '  -7.43%  '
$ws.Range('D42').Value = '''504.62'
$ws.Range('E42').Value = '  -5.70%  '
$ws.Range('D43').Value = '''24.58'
$ws.Range('E43').Value = '  -0.79%  '
$ws.Range('D44').Value = '''3.67'
$ws.Range('E44').Value = '  -3.51%  '
$ws.Range('E45').Value = '  -4.62%  '
$ws.Range('D46').Value = '''8.59'
$ws.Range('E46').Value = '  +5.49%  '
$ws.Range('D47').Value = '''0.0410'
$ws.Range('E47').Value = '  -7.15%  '
$ws.Range('D48').Value = '''1.65'
$ws.Range('E48').Value = '  +4.26%  '
$ws.Range('D49').Value = '''5.45'
$ws.Range('E49').Value = '  +2.40%  '
$ws.Range('D50').Value = '''53.19'
$ws.Range('E50').Value = '  +4.80%  '
$ws.Range('D51').Value = '''3.13'
$ws.Range('E51').Value = '  -11.07%  '
